$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 13
$ws.Cells.Item($row, 1).Value = 12
$ws.Cells.Item($row, 2).Value = "Monstack Developer"
$ws.Cells.Item($row, 3).Value = "Demoo"
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
